# Edit script: update Chirimoya - Macroferia Regional de Talca weekly price records
# Inserts 2 new observation rows at the top of the data block (rows 125-126),
# shifting existing observations down by two rows, and appends the final two
# displaced rows at the bottom (new rows 173-174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 125-172 with shifted values ---
# Row 125
$ws.Range("D125").Value = 45215
$ws.Range("L125").Value = 'Primera'
$ws.Range("M125").Value = 180
$ws.Range("N125").Value = 20000
$ws.Range("O125").Value = 20000
$ws.Range("P125").Value = 20000
$ws.Range("S125").Value = 2000
# Row 126
$ws.Range("D126").Value = 45215
$ws.Range("L126").Value = 'Segunda'
$ws.Range("M126").Value = 150
$ws.Range("N126").Value = 18000
$ws.Range("O126").Value = 18000
$ws.Range("P126").Value = 18000
$ws.Range("S126").Value = 1800
# Row 127
$ws.Range("D127").Value = 44421
$ws.Range("M127").Value = 30
$ws.Range("N127").Value = 35000
$ws.Range("O127").Value = 35000
$ws.Range("P127").Value = 35000
$ws.Range("S127").Value = 3500
# Row 128
$ws.Range("D128").Value = 44473
$ws.Range("L128").Value = 'Primera'
$ws.Range("M128").Value = 200
$ws.Range("N128").Value = 28000
$ws.Range("O128").Value = 28000
$ws.Range("P128").Value = 28000
$ws.Range("S128").Value = 2800
# Row 129
$ws.Range("L129").Value = 'Especial'
$ws.Range("M129").Value = 60
$ws.Range("N129").Value = 22000
$ws.Range("O129").Value = 22000
$ws.Range("P129").Value = 22000
$ws.Range("S129").Value = 2200
# Row 130
$ws.Range("D130").Value = 44841
$ws.Range("L130").Value = 'Extra (doble especial)'
$ws.Range("M130").Value = 50
$ws.Range("N130").Value = 25000
$ws.Range("O130").Value = 25000
$ws.Range("P130").Value = 25000
$ws.Range("S130").Value = 2500
# Row 131
$ws.Range("D131").Value = 44841
$ws.Range("M131").Value = 40
$ws.Range("N131").Value = 20000
$ws.Range("O131").Value = 20000
$ws.Range("P131").Value = 20000
$ws.Range("S131").Value = 2000
# Row 132
$ws.Range("D132").Value = 44459
$ws.Range("M132").Value = 100
# Row 133
$ws.Range("D133").Value = 44495
$ws.Range("L133").Value = 'Primera'
$ws.Range("N133").Value = 25000
$ws.Range("O133").Value = 25000
$ws.Range("P133").Value = 25000
$ws.Range("S133").Value = 2500
# Row 134
$ws.Range("D134").Value = 44474
$ws.Range("L134").Value = 'Especial'
$ws.Range("M134").Value = 150
# Row 135
$ws.Range("D135").Value = 45204
$ws.Range("M135").Value = 150
$ws.Range("N135").Value = 20000
$ws.Range("O135").Value = 20000
$ws.Range("P135").Value = 20000
$ws.Range("S135").Value = 2000
# Row 136
$ws.Range("D136").Value = 44446
$ws.Range("M136").Value = 200
$ws.Range("N136").Value = 30000
$ws.Range("O136").Value = 30000
$ws.Range("P136").Value = 30000
$ws.Range("S136").Value = 3000
# Row 137
$ws.Range("D137").Value = 45195
$ws.Range("M137").Value = 250
$ws.Range("N137").Value = 22000
$ws.Range("O137").Value = 22000
$ws.Range("P137").Value = 22000
$ws.Range("S137").Value = 2200
# Row 138
$ws.Range("D138").Value = 44494
$ws.Range("L138").Value = 'Primera'
$ws.Range("M138").Value = 150
$ws.Range("N138").Value = 25000
$ws.Range("O138").Value = 25000
$ws.Range("P138").Value = 25000
$ws.Range("S138").Value = 2500
# Row 139
$ws.Range("D139").Value = 44494
$ws.Range("L139").Value = 'Segunda'
$ws.Range("M139").Value = 50
$ws.Range("N139").Value = 23000
$ws.Range("O139").Value = 23000
$ws.Range("P139").Value = 23000
$ws.Range("S139").Value = 2300
# Row 140
$ws.Range("D140").Value = 45211
$ws.Range("M140").Value = 100
$ws.Range("N140").Value = 22000
$ws.Range("O140").Value = 22000
$ws.Range("P140").Value = 22000
$ws.Range("S140").Value = 2200
# Row 141
$ws.Range("D141").Value = 45211
$ws.Range("M141").Value = 230
$ws.Range("N141").Value = 20000
$ws.Range("O141").Value = 20000
$ws.Range("P141").Value = 20000
$ws.Range("S141").Value = 2000
# Row 142
$ws.Range("D142").Value = 44838
$ws.Range("L142").Value = 'Especial'
$ws.Range("M142").Value = 50
$ws.Range("N142").Value = 27000
$ws.Range("O142").Value = 27000
$ws.Range("P142").Value = 27000
$ws.Range("S142").Value = 2700
# Row 143
$ws.Range("D143").Value = 44838
$ws.Range("L143").Value = 'Primera'
$ws.Range("M143").Value = 60
$ws.Range("N143").Value = 24000
$ws.Range("O143").Value = 24000
$ws.Range("P143").Value = 24000
$ws.Range("S143").Value = 2400
# Row 144
$ws.Range("D144").Value = 45209
$ws.Range("M144").Value = 340
$ws.Range("N144").Value = 19000
$ws.Range("O144").Value = 19000
$ws.Range("P144").Value = 19000
$ws.Range("S144").Value = 1900
# Row 145
$ws.Range("D145").Value = 45209
$ws.Range("L145").Value = 'Segunda'
$ws.Range("M145").Value = 280
$ws.Range("N145").Value = 15000
$ws.Range("O145").Value = 15000
$ws.Range("P145").Value = 15000
$ws.Range("S145").Value = 1500
# Row 146
$ws.Range("D146").Value = 44799
$ws.Range("L146").Value = 'Primera'
$ws.Range("M146").Value = 50
$ws.Range("N146").Value = 30000
$ws.Range("O146").Value = 30000
$ws.Range("P146").Value = 30000
$ws.Range("S146").Value = 3000
# Row 147
$ws.Range("D147").Value = 44498
$ws.Range("M147").Value = 250
$ws.Range("N147").Value = 22000
$ws.Range("O147").Value = 23000
$ws.Range("P147").Value = 22600
$ws.Range("S147").Value = 2260
# Row 148
$ws.Range("L148").Value = 'Especial'
$ws.Range("M148").Value = 200
$ws.Range("N148").Value = 23000
$ws.Range("O148").Value = 23000
$ws.Range("P148").Value = 23000
$ws.Range("S148").Value = 2300
# Row 149
$ws.Range("D149").Value = 44859
$ws.Range("L149").Value = 'Primera'
$ws.Range("M149").Value = 150
$ws.Range("N149").Value = 20000
$ws.Range("O149").Value = 20000
$ws.Range("P149").Value = 20000
$ws.Range("S149").Value = 2000
# Row 150
$ws.Range("D150").Value = 44859
$ws.Range("L150").Value = 'Segunda'
$ws.Range("M150").Value = 120
$ws.Range("N150").Value = 18000
$ws.Range("O150").Value = 18000
$ws.Range("P150").Value = 18000
$ws.Range("S150").Value = 1800
# Row 151
$ws.Range("D151").Value = 44845
$ws.Range("L151").Value = 'Especial'
$ws.Range("M151").Value = 40
$ws.Range("N151").Value = 22000
$ws.Range("O151").Value = 22000
$ws.Range("P151").Value = 22000
$ws.Range("S151").Value = 2200
# Row 152
$ws.Range("D152").Value = 44845
$ws.Range("L152").Value = 'Primera'
$ws.Range("M152").Value = 40
# Row 153
$ws.Range("D153").Value = 44518
$ws.Range("M153").Value = 210
$ws.Range("N153").Value = 20000
$ws.Range("O153").Value = 20000
$ws.Range("P153").Value = 20000
$ws.Range("S153").Value = 2000
# Row 154
$ws.Range("D154").Value = 45190
$ws.Range("L154").Value = 'Segunda'
$ws.Range("M154").Value = 230
# Row 155
$ws.Range("D155").Value = 44438
$ws.Range("L155").Value = 'Primera'
$ws.Range("M155").Value = 100
$ws.Range("N155").Value = 30000
$ws.Range("O155").Value = 30000
$ws.Range("P155").Value = 30000
$ws.Range("S155").Value = 3000
# Row 156
$ws.Range("D156").Value = 45212
$ws.Range("L156").Value = 'Primera'
$ws.Range("M156").Value = 180
$ws.Range("N156").Value = 20000
$ws.Range("O156").Value = 20000
$ws.Range("P156").Value = 20000
$ws.Range("S156").Value = 2000
# Row 157
$ws.Range("D157").Value = 45212
$ws.Range("L157").Value = 'Segunda'
$ws.Range("M157").Value = 150
$ws.Range("N157").Value = 18000
$ws.Range("O157").Value = 18000
$ws.Range("P157").Value = 18000
$ws.Range("S157").Value = 1800
# Row 158
$ws.Range("L158").Value = 'Especial'
$ws.Range("M158").Value = 100
$ws.Range("N158").Value = 30000
$ws.Range("O158").Value = 30000
$ws.Range("P158").Value = 30000
$ws.Range("S158").Value = 3000
# Row 159
$ws.Range("D159").Value = 44813
$ws.Range("L159").Value = 'Primera'
$ws.Range("M159").Value = 130
$ws.Range("N159").Value = 25000
$ws.Range("O159").Value = 25000
$ws.Range("P159").Value = 25000
$ws.Range("S159").Value = 2500
# Row 160
$ws.Range("D160").Value = 44813
$ws.Range("L160").Value = 'Segunda'
$ws.Range("M160").Value = 60
$ws.Range("N160").Value = 20000
$ws.Range("O160").Value = 20000
$ws.Range("P160").Value = 20000
$ws.Range("S160").Value = 2000
# Row 161
$ws.Range("D161").Value = 44434
$ws.Range("L161").Value = 'Especial'
$ws.Range("M161").Value = 60
# Row 162
$ws.Range("D162").Value = 44511
$ws.Range("L162").Value = 'Primera'
$ws.Range("M162").Value = 200
$ws.Range("N162").Value = 25000
$ws.Range("O162").Value = 25000
$ws.Range("P162").Value = 25000
$ws.Range("S162").Value = 2500
# Row 163
$ws.Range("D163").Value = 44441
$ws.Range("L163").Value = 'Primera'
$ws.Range("M163").Value = 150
# Row 164
$ws.Range("D164").Value = 44432
$ws.Range("M164").Value = 70
$ws.Range("N164").Value = 30000
$ws.Range("O164").Value = 30000
$ws.Range("P164").Value = 30000
$ws.Range("S164").Value = 3000
# Row 165
$ws.Range("D165").Value = 44803
$ws.Range("L165").Value = 'Especial'
$ws.Range("M165").Value = 40
$ws.Range("N165").Value = 30000
$ws.Range("O165").Value = 30000
$ws.Range("P165").Value = 30000
$ws.Range("S165").Value = 3000
# Row 166
$ws.Range("L166").Value = 'Especial'
$ws.Range("M166").Value = 130
$ws.Range("N166").Value = 23000
$ws.Range("O166").Value = 23000
$ws.Range("P166").Value = 23000
$ws.Range("S166").Value = 2300
# Row 167
$ws.Range("D167").Value = 44900
$ws.Range("L167").Value = 'Extra (doble especial)'
$ws.Range("M167").Value = 70
$ws.Range("N167").Value = 25000
$ws.Range("O167").Value = 25000
$ws.Range("P167").Value = 25000
$ws.Range("S167").Value = 2500
# Row 168
$ws.Range("D168").Value = 44900
$ws.Range("M168").Value = 180
$ws.Range("N168").Value = 20000
$ws.Range("O168").Value = 20000
$ws.Range("P168").Value = 20000
$ws.Range("S168").Value = 2000
# Row 169
$ws.Range("L169").Value = 'Especial'
$ws.Range("M169").Value = 150
$ws.Range("N169").Value = 28000
$ws.Range("O169").Value = 28000
$ws.Range("P169").Value = 28000
$ws.Range("S169").Value = 2800
# Row 170
$ws.Range("D170").Value = 44832
$ws.Range("L170").Value = 'Primera'
$ws.Range("M170").Value = 200
$ws.Range("N170").Value = 25000
$ws.Range("O170").Value = 25000
$ws.Range("P170").Value = 25000
$ws.Range("S170").Value = 2500
# Row 171
$ws.Range("D171").Value = 44832
$ws.Range("L171").Value = 'Segunda'
$ws.Range("M171").Value = 180
$ws.Range("N171").Value = 20000
$ws.Range("O171").Value = 20000
$ws.Range("P171").Value = 20000
$ws.Range("S171").Value = 2000
# Row 172
$ws.Range("D172").Value = 44455
$ws.Range("L172").Value = 'Especial'
$ws.Range("M172").Value = 150
$ws.Range("N172").Value = 30000
$ws.Range("O172").Value = 30000
$ws.Range("P172").Value = 30000
$ws.Range("S172").Value = 3000

# --- Add new rows 173 and 174 (same template as row 172, new D/L/M/N/O/P/S) ---
# Row 173
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = 'Macroferia Regional de Talca'
$ws.Range("C173").Value = 'Maule'
$ws.Range("D173").NumberFormat = $ws.Range("D124").NumberFormat
$ws.Range("D173").Value = 44809
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 'Fruta'
$ws.Range("G173").Value = 100107
$ws.Range("H173").Value = 'Otros'
$ws.Range("I173").Value = 100107002
$ws.Range("J173").Value = 'Chirimoya'
$ws.Range("K173").Value = 'Cultivar IV Región'
$ws.Range("L173").Value = 'Especial'
$ws.Range("M173").Value = 150
$ws.Range("N173").Value = 28000
$ws.Range("O173").Value = 28000
$ws.Range("P173").Value = 28000
$ws.Range("Q173").Value = '$/bandeja 10 kilos'
$ws.Range("R173").Value = 'Provincia de Limarí'
$ws.Range("S173").Value = 2800
$ws.Range("T173").Value = 10

# Row 174
$ws.Range("A174").Value = 5
$ws.Range("B174").Value = 'Macroferia Regional de Talca'
$ws.Range("C174").Value = 'Maule'
$ws.Range("D174").NumberFormat = $ws.Range("D124").NumberFormat
$ws.Range("D174").Value = 44809
$ws.Range("E174").Value = 7
$ws.Range("F174").Value = 'Fruta'
$ws.Range("G174").Value = 100107
$ws.Range("H174").Value = 'Otros'
$ws.Range("I174").Value = 100107002
$ws.Range("J174").Value = 'Chirimoya'
$ws.Range("K174").Value = 'Cultivar IV Región'
$ws.Range("L174").Value = 'Primera'
$ws.Range("M174").Value = 100
$ws.Range("N174").Value = 25000
$ws.Range("O174").Value = 25000
$ws.Range("P174").Value = 25000
$ws.Range("Q174").Value = '$/bandeja 10 kilos'
$ws.Range("R174").Value = 'Provincia de Limarí'
$ws.Range("S174").Value = 2500
$ws.Range("T174").Value = 10
